$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.802.84"
$ws.Range("E2").Value = "  -0.74%  "

$ws.Range("D3").Value = "3.430.05"
$ws.Range("E3").Value = "  -1.91%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.05%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.69%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  -1.28%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.59"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.08%  "

$ws.Range("E10").Value = "  +1.12%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.383"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.64%  "

$ws.Range("D12").Value = "4.012.90"
$ws.Range("E12").Value = "  -1.95%  "

$ws.Range("E13").Value = "  -0.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000177"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.65%  "

$ws.Range("D15").Value = "3.425.39"
$ws.Range("E15").Value = "  -2.07%  "

$ws.Range("D16").Value = "63.821.70"
$ws.Range("E16").Value = "  -0.84%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.38%  "

$ws.Range("E18").Value = "  +0.37%  "

$ws.Range("E19").Value = "  -1.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.25%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "385.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.99%  "

$ws.Range("E22").Value = "  -1.10%  "

$ws.Range("D23").Value = "3.568.67"
$ws.Range("E23").Value = "  -1.87%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.85%  "

$ws.Range("E25").Value = "  +0.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000110"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.38%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.996"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.52%  "

$ws.Range("E28").Value = "  -1.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.85%  "

$ws.Range("E30").Value = "  -3.24%  "

$ws.Range("E31").Value = "  +1.86%  "

$ws.Range("E32").Value = "  -3.40%  "

$ws.Range("D33").Value = "3.459.70"
$ws.Range("E33").Value = "  -1.71%  "

$ws.Range("E34").Value = "  -0.08%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "22.93"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.15%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.20"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.19%  "

$ws.Range("E37").Value = "  -1.29%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "164.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.77%  "

$ws.Range("E39").Value = "  -2.14%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0775"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.69%  "

$ws.Range("E41").Value = "  -2.37%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.10%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.00%  "

$ws.Range("E44").Value = "  -0.95%  "

$ws.Range("E45").Value = "  -2.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.54"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.40%  "

$ws.Range("E47").Value = "  -4.15%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.10%  "

$ws.Range("E49").Value = "  +0.90%  "

$ws.Range("D50").Value = "2.289.45"
$ws.Range("E50").Value = "  -7.15%  "

$ws.Range("E51").Value = "  -2.09%  "
